$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "Sunday 22.4.18"
$ws.Range("A5").Value = "Sunday 15.4.18"
$ws.Range("B5").Value = "1000-1600"
$ws.Range("C5").Value = 6
$ws.Range("D5").Value = "Getting the bastard to work"
$ws.Range("B6").Value = "1000-1230"
$ws.Range("C6").Value = 2.5
$ws.Range("D6").Value = "Trying to connect to device"

$ws.Range("A8").Select()
